$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the resolved case 5894 (ALBARELLOS AV. 3100) - row 3.
#    This shifts all subsequent rows up by one position.
# ------------------------------------------------------------------
$ws.Rows(3).Delete()

# ------------------------------------------------------------------
# 2. Apply field updates to existing cases (rows after the shift).
# ------------------------------------------------------------------
# Case 4238 (now row 3): Observaciones updated
    $ws.Cells.Item(3, 8).Value = 'cable en panza'

# Case 6557 (now row 5): OT updated
    $ws.Cells.Item(5, 5).Value = 'Pendiente ADM'

# Case 6569 (now row 7): OT updated
    $ws.Cells.Item(7, 5).Value = 'ICD30326286 '

# ------------------------------------------------------------------
# 3. Append the 4 new cases reported at the bottom of the sheet.
# ------------------------------------------------------------------
# Row 8: case 6583 - MEXICO 1942
    $ws.Cells.Item(8, 1).Value = '''6583'
    $ws.Cells.Item(8, 2).Value = '''8/7/2025'
    $ws.Cells.Item(8, 3).Value = 'MEXICO 1942'
    $ws.Cells.Item(8, 4).Value = 3
    $ws.Cells.Item(8, 5).Value = 'Pendiente ADM'
    $ws.Cells.Item(8, 6).Value = 'Optical Power'
    $ws.Cells.Item(8, 7).Value = 'Pendiente'
    $ws.Cells.Item(8, 8).Value = 'Cables en panza cortados y suelto'
    $ws.Cells.Item(8, 9).Value = 1
    $ws.Cells.Item(8, 10).Value = '{"direccionesNormalizadas": [{"altura": 1942, "cod_calle": 13076, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.393874", "y": "-34.615880"}, "direccion": "MEXICO 1942, CABA", "nombre_calle": "MEXICO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
    $ws.Cells.Item(8, 11).Value = -58.393874
    $ws.Cells.Item(8, 12).Value = -34.61588
    $ws.Cells.Item(8, 13).Value = 'Almagro'
    $ws.Cells.Item(8, 14).Value = 'Capital Sur'

# Row 9: case 6100 - DE LOS CONSTITUYENTES AV. 5552
    $ws.Cells.Item(9, 1).Value = '''6100'
    $ws.Cells.Item(9, 2).Value = '''8/8/2025'
    $ws.Cells.Item(9, 3).Value = 'DE LOS CONSTITUYENTES AV. 5552'
    $ws.Cells.Item(9, 4).Value = 12
    $ws.Cells.Item(9, 5).Value = 'Pendiente ADM'
    $ws.Cells.Item(9, 6).Value = 'Optical Power'
    $ws.Cells.Item(9, 7).Value = 'Pendiente'
    $ws.Cells.Item(9, 8).Value = 'Cable en panza y cables cortados'
    $ws.Cells.Item(9, 9).Value = 1
    $ws.Cells.Item(9, 10).Value = '{"direccionesNormalizadas": [{"altura": 5552, "cod_calle": 4043, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.501174", "y": "-34.575005"}, "direccion": "DE LOS CONSTITUYENTES AV. 5552, CABA", "nombre_calle": "DE LOS CONSTITUYENTES AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
    $ws.Cells.Item(9, 11).Value = -58.501174
    $ws.Cells.Item(9, 12).Value = -34.575005
    $ws.Cells.Item(9, 13).Value = 'Paternal'
    $ws.Cells.Item(9, 14).Value = 'Capital Norte'

# Row 10: case 6154 - PATAGONES 2728
    $ws.Cells.Item(10, 1).Value = '''6154'
    $ws.Cells.Item(10, 2).Value = '''8/7/2025'
    $ws.Cells.Item(10, 3).Value = 'PATAGONES 2728'
    $ws.Cells.Item(10, 4).Value = 4
    $ws.Cells.Item(10, 5).Value = 'Pendiente ADM'
    $ws.Cells.Item(10, 6).Value = 'Optical Power'
    $ws.Cells.Item(10, 7).Value = 'Pendiente'
    $ws.Cells.Item(10, 8).Value = 'Cable en panza cables cortados cables sueltos'
    $ws.Cells.Item(10, 9).Value = 1
    $ws.Cells.Item(10, 10).Value = '{"direccionesNormalizadas": [{"altura": 2728, "cod_calle": 17032, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.400609", "y": "-34.637438"}, "direccion": "PATAGONES 2728, CABA", "nombre_calle": "PATAGONES", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
    $ws.Cells.Item(10, 11).Value = -58.400609
    $ws.Cells.Item(10, 12).Value = -34.637438
    $ws.Cells.Item(10, 13).Value = 'San Telmo'
    $ws.Cells.Item(10, 14).Value = 'Capital Sur'

# Row 11: case 6265 - BROWN, ALTE. AV. 881
    $ws.Cells.Item(11, 1).Value = '''6265'
    $ws.Cells.Item(11, 2).Value = '''8/7/2025'
    $ws.Cells.Item(11, 3).Value = 'BROWN, ALTE. AV. 881'
    $ws.Cells.Item(11, 4).Value = 4
    $ws.Cells.Item(11, 5).Value = 'Pendiente ADM'
    $ws.Cells.Item(11, 6).Value = 'Optical Power'
    $ws.Cells.Item(11, 7).Value = 'Pendiente'
    $ws.Cells.Item(11, 8).Value = 'tendido a baja altura'
    $ws.Cells.Item(11, 9).Value = 1
    $ws.Cells.Item(11, 10).Value = '{"direccionesNormalizadas": [{"altura": 881, "cod_calle": 2115, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.360551", "y": "-34.632684"}, "direccion": "BROWN, ALTE. AV. 881, CABA", "nombre_calle": "BROWN, ALTE. AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
    $ws.Cells.Item(11, 11).Value = -58.360551
    $ws.Cells.Item(11, 12).Value = -34.632684
    $ws.Cells.Item(11, 13).Value = 'San Telmo'
    $ws.Cells.Item(11, 14).Value = 'Capital Sur'

